$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    "sweat pants 6x",
    "sweat pants womans petites",
    "sweat tights",
    "sweat tights women",
    "sweaty betty",
    "swelling machine",
    "swimmers body wash",
    "swimming compression shorts",
    "swimming leggings for women",
    "tall leggings",
    "tall leggings for women",
    "tall length pants women",
    "tall opaque black tights",
    "tall running pants",
    "tall tights for women",
    "tall tights for women long",
    "tall tights women",
    "tall womens tights",
    "tape clothes to skin",
    "tape for knee support",
    "tape for muscle support",
    "tape for muscles pain",
    "tape pants",
    "tape support",
    "tape to hold dress in place",
    "tape underwear",
    "tape zebra",
    "target clothes",
    "target clothes for women",
    "target leggings for women",
    "target pants for women",
    "target sports",
    "team leggings for women",
    "team warm up pants",
    "tear away basketball pants",
    "tear away leggings",
    "tear away pants for women",
    "tear away track pants women",
    "tears run rings",
    "tech capris",
    "tech gear leggings",
    "tech gear yoga pants",
    "tech pants for women",
    "tech suits for women",
    "tech suits swimming women",
    "tech wash",
    "technical pants",
    "technical support",
    "technical tights",
    "teen activity tracker",
    "teen black leggings",
    "teen tights",
    "teen weight training",
    "tengo womens ripped hole bandage pants denim shorts",
    "tennis alone",
    "tennis capri skirt",
    "tennis capri women",
    "tennis capris",
    "tennis capris with pockets",
    "tennis conditioning",
    "tennis gear for women",
    "tennis knee band",
    "tennis leggings",
    "tennis long pants",
    "tennis pants",
    "tennis pants for women with pockets",
    "tennis pants women",
    "tennis pants womens",
    "tennis skirt leggings",
    "tennis skirt with leggings",
    "tennis skirt with leggings women",
    "tennis skirts long",
    "tennis skirts with leggings",
    "tennis tights",
    "tennis tommy for women",
    "tennis warm up pants",
    "tennis warmups for men",
    "tesla compression pants women",
    "tganksgiving leggings",
    "thailand clothes for women",
    "thailand pants",
    "thailand pants men",
    "thailand women clothing",
    "the amazing race 22",
    "the amazing race 25",
    "the best leggings",
    "the bread club",
    "the cw live",
    "the endurance club",
    "the endurance model",
    "the essential brace",
    "the eve exo",
    "the flash compression pants",
    "the flash compression shirt",
    "the flash sleep pants for men",
    "the flash sweat pants",
    "the flash underwear women",
    "the flash workout clothes",
    "the flash workout shirt",
    "the game of things nasty"
)

for ($i = 0; $i -lt $values.Count; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

